$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 5351.066666666667, 5697, 4691, 0.1155822277069092),
    @(1, 5241.366666666667, 5584, 4762, 0.1173498868942261),
    @(2, 5180.066666666667, 5603, 4199, 0.1183500607808431),
    @(3, 5670.966666666666, 5963, 5325, 0.11767737865448),
    @(4, 5042.566666666667, 5358, 4432, 0.1222634712855021),
    @(5, 5010.633333333333, 5413, 4378, 0.1175279140472412),
    @(6, 5529.1, 5984, 4872, 0.120494016011556),
    @(7, 5383.833333333333, 5789, 4916, 0.121275266011556),
    @(8, 5328.633333333333, 5577, 4668, 0.1213013569513957),
    @(9, 5333.666666666667, 5927, 4877, 0.11748259862264)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
